# Apply timetable updates: Section_A, Section_B cell changes, and add a
# new "Course_Summary" sheet describing each course code.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Section_A (sheet1) updates
# ---------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("C2").Value = "Free"
$wsA.Range("E2").Value = "Free"

$wsA.Range("B3").Value = "CS309 (Tutorial)"

$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "CS303"
$wsA.Range("E5").Value = "CS304"
$wsA.Range("F5").Value = "CS309"

$wsA.Range("B6").Value = "CS304 (Tutorial)"
$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "CS309"
$wsA.Range("F6").Value = "CS304"

$wsA.Range("B7").Value = "CS304"
$wsA.Range("D7").Value = "Free"
$wsA.Range("E7").Value = "CS303 (Tutorial)"
$wsA.Range("F7").Value = "Free"

# ---------------------------------------------------------------------
# Section_B (sheet2) updates
# ---------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("C2").Value = "CS304 (Tutorial)"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "CS309 (Tutorial)"

$wsB.Range("B3").Value = "CS303 (Tutorial)"
$wsB.Range("C3").Value = "Free"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "CS304"

$wsB.Range("C5").Value = "CS309"
$wsB.Range("D5").Value = "CS303"
$wsB.Range("E5").Value = "CS304"
$wsB.Range("F5").Value = "CS309"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "CS304"
$wsB.Range("D6").Value = "Free"
$wsB.Range("E6").Value = "CS303"
$wsB.Range("F6").Value = "Free"

$wsB.Range("B7").Value = "CS303"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "CS309"
$wsB.Range("E7").Value = "Free"
$wsB.Range("F7").Value = "Free"

# ---------------------------------------------------------------------
# New "Course_Summary" sheet, placed after Section_B
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$wsSummary = $wb.Worksheets.Add($null, $lastSheet)
$wsSummary.Name = "Course_Summary"

$headers = @("Course Code","Course Name","Course Type","LTPSC","Credits","Instructor")
for ($col = 0; $col -lt $headers.Length; $col++) {
    $wsSummary.Cells.Item(1, $col + 1).Value = $headers[$col]
}

# Match the bold/centered/bordered header formatting used on the other sheets.
$wsA.Range("B1").Copy()
$wsSummary.Range("A1:F1").PasteSpecial(-4122)

$courses = @(
    @("CS309", "Statistics for Computer Science", "Core",     "3-1-0-0-4", 4, "Dr. Sunil C K"),
    @("CS303", "Computer Networks",                "Core",     "3-1-2-0-6", 6, "Dr. Animesh Roy"),
    @("CS304", "Artificial Intelligence",           "Core",     "3-1-0-0-4", 4, "Dr. Krishnendu"),
    @("CS461", "Parallel computing",                "Elective", "4-0-0-0-4", 4, "Dr. Pramod")
)

$rowIndex = 2
foreach ($course in $courses) {
    for ($col = 0; $col -lt $course.Length; $col++) {
        $wsSummary.Cells.Item($rowIndex, $col + 1).Value = $course[$col]
    }
    $rowIndex++
}

# Keep the originally-active sheet selected (adding a sheet makes it active).
$wsA.Activate()
[void]$wsA.Range("A1").Select()
